# Summary.xlsx edit: "Fixed cigar tables in results"
#
# The author inserted a new data column (a new "LM00" sample) into the
# comparison table. In Excel terms this is:
#   1. Insert a new column at D (pushes old D..F -> E..G for every row).
#   2. Fill in the two new data points for the new column (rows 2 & 3).
#   3. Fix up the header row, which the author rearranged by hand:
#        the old leftmost header ("Next", originally in A1) got shifted
#        one column to the right (to B1) rather than staying at A1,
#        and the new header label "LM00" was typed into the freshly
#        inserted column's header cell.
#   4. A couple of formula ranges were re-entered as block fills, which
#      Excel groups into shared formulas - replicate that with ranged
#      Formula assignment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before D; this shifts D:F -> E:G on every row
#        (rows 2-23), which matches the data/formula shift seen for all
#        rows except the header row.
$ws.Columns("D:D").Insert()

# --- 2. New data values for the inserted column (only rows 2 and 3 have
#        figures for the new "LM00" sample).
$ws.Range("D2").Value2 = 12.166
$ws.Range("D3").Value2 = 144.9974

# --- 3. Rebuild the header row (row 1) to its final arrangement.
#        After the column insert it reads: A1=Next B1=Marlbro C1=LM
#        D1=(blank) E1=Galouis F1=Canadian
#        Final layout needed: B1=Next C1=Marlbro D1=LM E1=LM00 F1=Galouis
#        G1=Canadian (A1 cleared).
# (Note: use Value2, not Value, to read/write cell contents here.)
$ws.Range("G1").Value2 = $ws.Range("F1").Value2
$ws.Range("F1").Value2 = $ws.Range("E1").Value2
$ws.Range("E1").Value2 = "LM00"
$ws.Range("D1").Value2 = $ws.Range("C1").Value2
$ws.Range("C1").Value2 = $ws.Range("B1").Value2
$ws.Range("B1").Value2 = $ws.Range("A1").Value2
$ws.Range("A1").ClearContents()

# --- 4. Re-enter a couple of formula blocks in one shot so they collapse
#        into shared formulas, matching how Excel compacts repeated
#        formulas entered/filled together.
$ws.Range("B16:B17").Formula = "=C7/C$11 * 100"
$ws.Range("B20:B21").Formula = "=C2/C$11 * 100"
$ws.Range("C22").Formula = "=SUM(C14:C21)"

# Make the same cell active/selected as in the final saved workbook.
$ws.Range("J6").Select()
